$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.076.98"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "'1.645.03"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.76%  "
$ws.Range("D5").Value = "'215.99"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'0.507"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "'0.0637"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").Value = "'19.46"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "'0.0796"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.663.46"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.25"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "'63.39"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "'0.0₃0759"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "'26.080.13"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "'194.17"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").Value = "'9.75"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("D24").Value = "'144.33"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "'1.78"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "'6.86"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'15.53"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'1.25"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "'0.0488"
$ws.Range("E30").Value = "  -2.81%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "'0.900"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'1.129.92"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'0.532"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.798"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "'98.53"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Value = "'0.0₆0113"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D44").Value = "'56.22"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").Value = "'0.0522"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "'7.70"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").Value = "'0.417"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "'0.0942"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("E51").Value = "  +1.21%  "
